$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = 2.25048404231296
$ws.Range("D8").Value = 0.145502777593941
$ws.Range("E8").Value = 0.588461538461538
$ws.Range("F8").Value = 0.672727272727273
$ws.Range("G8").Value = 0.989295231124154
